# Fill in "Special Attack 2" column (E/F) for the remaining Good Guys
# (Archer, Druid, Necromancer, Warrior, Wizard) that were missing it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "RapidShot"
$ws.Range("F2").Value = "y"

$ws.Range("E4").Value = "SummonWolf"
$ws.Range("F4").Value = "y"

$ws.Range("E6").Value = "SummonSkeleton"
$ws.Range("F6").Value = "y"

$ws.Range("E8").Value = "WhirlWindOfDeath"
$ws.Range("F8").Value = "y"

$ws.Range("E9").Value = "AvadaKedavra"
$ws.Range("F9").Value = "y"

# Update the view so the selected cell is F9 and scroll position resets
# to the top of the sheet (matches the saved window state in the diff).
$ws.Range("F9").Select()
